$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H42").Value = 274.75
$ws.Range("I42").Value = 100
$ws.Range("J42").Value = 333
$ws.Range("K42").Value = 300
$ws.Range("L42").Value = 999
$ws.Range("M42").Value = -70
$ws.Range("N42").Value = -1459

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 14803.6
$ws.Range("J43").Value = 9936.916999999999
$ws.Range("L43").Value = 9936.916999999999
$ws.Range("N43").Value = -10074.917

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 5213.9
$ws.Range("I92").Value = 5987.0586
$ws.Range("K92").Value = 5987.0586
$ws.Range("M92").Value = -4739.0586

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 11178886
$ws.Range("I132").Value = 11529882
$ws.Range("K132").Value = 34589646
$ws.Range("M132").Value = -34587116

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 1395.6923
$ws.Range("I135").Value = 914.4
$ws.Range("K135").Value = 8229.6
$ws.Range("M135").Value = -5694.6

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 4203.5405
$ws.Range("I138").Value = 1931.3334
$ws.Range("J138").Value = 5294.2
$ws.Range("K138").Value = 5794.0002
$ws.Range("L138").Value = 15882.6
$ws.Range("M138").Value = -654.0002000000004
$ws.Range("N138").Value = -26162.6

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 200966.67
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 200966.67
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 602900.01
$ws.Range("M141").Value = ""
$ws.Range("N141").Value = -613260.01

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3903.6128
$ws.Range("I32").Value = 3490.0688
$ws.Range("K32").Value = 3490.0688
$ws.Range("M32").Value = -3203.0688

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H37").Value = 950
$ws.Range("I37").Value = 950
$ws.Range("K37").Value = 950
$ws.Range("M37").Value = -677

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H55").Value = 13415.333
$ws.Range("I55").Value = 248
$ws.Range("K55").Value = 248
$ws.Range("M55").Value = 67

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 3706.8572
$ws.Range("I132").Value = 3602.5
$ws.Range("K132").Value = 10807.5
$ws.Range("M132").Value = -8277.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 4456.0464
$ws.Range("I20").Value = 4349.759
$ws.Range("K20").Value = 4349.759
$ws.Range("M20").Value = -4102.759

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 13605.448
$ws.Range("I86").Value = 9726.571
$ws.Range("K86").Value = 9726.571
$ws.Range("M86").Value = -8603.571

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 13605.448
$ws.Range("I89").Value = 9726.571
$ws.Range("K89").Value = 48632.855
$ws.Range("M89").Value = -43016.855

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1518.2593
$ws.Range("I94").Value = 1190.1428
$ws.Range("J94").Value = 2666.6667
$ws.Range("K94").Value = 1190.1428
$ws.Range("L94").Value = 2666.6667
$ws.Range("M94").Value = -739.1428000000001
$ws.Range("N94").Value = -3568.6667

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2526.4517
$ws.Range("I105").Value = 2534.1738
$ws.Range("K105").Value = 2534.1738
$ws.Range("M105").Value = -787.1738

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 1565.6666
$ws.Range("I94").Value = 1726
$ws.Range("K94").Value = 1726
$ws.Range("M94").Value = -1275

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 1947.2858
$ws.Range("I105").Value = 1842
$ws.Range("J105").Value = 2333.3333
$ws.Range("K105").Value = 1842
$ws.Range("L105").Value = 2333.3333
$ws.Range("M105").Value = -95
$ws.Range("N105").Value = -5827.3333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 14282.393
$ws.Range("I134").Value = 6430.391
$ws.Range("J134").Value = 50401.6
$ws.Range("K134").Value = 19291.173
$ws.Range("L134").Value = 151204.8
$ws.Range("M134").Value = -16756.173
$ws.Range("N134").Value = -156274.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 4936789.5
$ws.Range("I4").Value = 5334771.5
$ws.Range("K4").Value = 16004314.5
$ws.Range("M4").Value = -16004202.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H29").Value = 755.5714
$ws.Range("I29").Value = 357.8
$ws.Range("J29").Value = 1750
$ws.Range("K29").Value = 1073.4
$ws.Range("L29").Value = 5250
$ws.Range("M29").Value = -796.4000000000001
$ws.Range("N29").Value = -5804

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H141").Value = 171584.83
$ws.Range("I141").Value = 5895.2
$ws.Range("K141").Value = 17685.6
$ws.Range("M141").Value = -12505.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 9709.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 9299.388999999999
$ws.Range("I70").Value = 10672
$ws.Range("J70").Value = 7142.4287
$ws.Range("K70").Value = 10672
$ws.Range("L70").Value = 7142.4287
$ws.Range("M70").Value = -10402
$ws.Range("N70").Value = -7682.4287

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 9299.388999999999
$ws.Range("I73").Value = 10672
$ws.Range("J73").Value = 7142.4287
$ws.Range("K73").Value = 10672
$ws.Range("L73").Value = 7142.4287
$ws.Range("M73").Value = -9736
$ws.Range("N73").Value = -9014.4287

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1456.6666
$ws.Range("I97").Value = 1373.1538
$ws.Range("K97").Value = 1373.1538
$ws.Range("M97").Value = -877.1538

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 17650.3
$ws.Range("I7").Value = 18611.445
$ws.Range("K7").Value = 18611.445
$ws.Range("M7").Value = -18499.445

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 7744.3335
$ws.Range("I22").Value = 2316.6667
$ws.Range("J22").Value = 18599.666
$ws.Range("K22").Value = 2316.6667
$ws.Range("L22").Value = 18599.666
$ws.Range("M22").Value = -2021.6667
$ws.Range("N22").Value = -19189.666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 7744.3335
$ws.Range("I27").Value = 2316.6667
$ws.Range("J27").Value = 18599.666
$ws.Range("K27").Value = 2316.6667
$ws.Range("L27").Value = 18599.666
$ws.Range("M27").Value = -2209.6667
$ws.Range("N27").Value = -18813.666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4528.5713
$ws.Range("I40").Value = 3950
$ws.Range("J40").Value = 8000
$ws.Range("K40").Value = 3950
$ws.Range("L40").Value = 8000
$ws.Range("M40").Value = -3814
$ws.Range("N40").Value = -8272

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 246.07692
$ws.Range("I55").Value = 224.88889
$ws.Range("J55").Value = 293.75
$ws.Range("K55").Value = 224.88889
$ws.Range("L55").Value = 293.75
$ws.Range("M55").Value = -51.88889
$ws.Range("N55").Value = -639.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 6418.381
$ws.Range("I68").Value = 3778.9333
$ws.Range("J68").Value = 13017
$ws.Range("K68").Value = 3778.9333
$ws.Range("L68").Value = 13017
$ws.Range("M68").Value = -3029.9333
$ws.Range("N68").Value = -14515

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 6418.381
$ws.Range("I71").Value = 3778.9333
$ws.Range("J71").Value = 13017
$ws.Range("K71").Value = 18894.6665
$ws.Range("L71").Value = 65085
$ws.Range("M71").Value = -15150.6665
$ws.Range("N71").Value = -72573

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 657164.5
$ws.Range("I93").Value = 928199.75
$ws.Range("J93").Value = 6679.8
$ws.Range("K93").Value = 928199.75
$ws.Range("L93").Value = 6679.8
$ws.Range("M93").Value = -926951.75
$ws.Range("N93").Value = -9175.799999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 2528.5334
$ws.Range("I100").Value = 2535.6667
$ws.Range("K100").Value = 2535.6667
$ws.Range("M100").Value = -1994.6667

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 17650.3
$ws.Range("I126").Value = 18611.445
$ws.Range("K126").Value = 55834.335
$ws.Range("M126").Value = -53364.335

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3979.353
$ws.Range("I132").Value = 3300.5652
$ws.Range("J132").Value = 5398.636
$ws.Range("K132").Value = 9901.695599999999
$ws.Range("L132").Value = 16195.908
$ws.Range("M132").Value = -7371.695599999999
$ws.Range("N132").Value = -21255.908

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 4911.3213
$ws.Range("J136").Value = 4718.75
$ws.Range("L136").Value = 14156.25
$ws.Range("N136").Value = -19256.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3657.1765
$ws.Range("I132").Value = 3013
$ws.Range("J132").Value = 6663.3335
$ws.Range("K132").Value = 9039
$ws.Range("L132").Value = 19990.0005
$ws.Range("M132").Value = -6509
$ws.Range("N132").Value = -25050.0005
